# Apply the "full column transfer working with volumes" edit to the master sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master")

# Clear the stray customer_primer values ("3R"/"3F") that were left over in E3/E5.
$ws.Range("E3").ClearContents()
$ws.Range("E5").ClearContents()

# Row 12 (source_type = "strip") now has a real bcl_primer selected, driving the
# downstream H12/I12/K12 formulas (primer lookup, volume calc, remaining vol calc).
$ws.Range("F12").Value = "M13rev(-29)"

# Rows 26-33 (source_type = "plate") are now filled in as a full 8-well column
# transfer: column A of the source plate (A01..H01), one well per row.
# NOTE: column C is intentionally written before column B on each row so that the
# dependent cross-sheet formulas (source_plate!/source_strip!/source_tube!/all_sources!)
# are marked dirty and recalculate correctly.
$ws.Range("C26").Value = "A01"
$ws.Range("B26").Value = "plate"

$ws.Range("C27").Value = "B01"
$ws.Range("B27").Value = "plate"

$ws.Range("C28").Value = "C01"
$ws.Range("B28").Value = "plate"

$ws.Range("C29").Value = "D01"
$ws.Range("B29").Value = "plate"

$ws.Range("C30").Value = "E01"
$ws.Range("B30").Value = "plate"

$ws.Range("C31").Value = "F01"
$ws.Range("B31").Value = "plate"

$ws.Range("C32").Value = "G01"
$ws.Range("B32").Value = "plate"

$ws.Range("C33").Value = "H01"
$ws.Range("B33").Value = "plate"

# Reflect where the user was last working (selection/active-cell bookkeeping).
$allSources = $wb.Worksheets.Item("all_sources")
$allSources.Activate()
$allSources.Range("L23").Select() | Out-Null

$ws.Activate()
$ws.Range("E30").Select() | Out-Null

$wb.Save()
